# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.224.11'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '2.268.14'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D5").Value = '306.92'
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").Value = '97.01'
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("D7").Value = '0.526'
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("D10").Value = '35.26'
$ws.Range("E10").Value = '  +3.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '6.89'
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("D14").Value = '2.620.60'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = '14.74'
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").Value = '2.271.40'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '0.794'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '42.098.41'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.40'
$ws.Range("D20").Value = '0.0₃0906'
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = '6.03'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = '68.15'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.30'
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("D24").Value = '2.57'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '23.59'
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("E28").Value = '  +5.03%  '
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").Value = '2.11'
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("D31").Value = '161.78'
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D35").Value = '0.0737'
$ws.Range("E35").Value = '  -2.09%  '
$ws.Range("D36").Value = '17.23'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("E37").Value = '  -0.47%  '
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("D41").Value = '4.06'
$ws.Range("E41").Value = '  -4.02%  '
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("D43").Value = '1.947.81'
$ws.Range("E43").Value = '  -3.55%  '
$ws.Range("D44").Value = '18.89'
$ws.Range("E44").Value = '  -3.82%  '
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("D46").Value = '9.95'
$ws.Range("E46").Value = '  -2.82%  '
$ws.Range("D47").Value = '2.88'
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.60'
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").Value = '71.92'
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").Value = '92.14'
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("E51").Value = '  -1.60%  '
